# Update cached FFXIV market-board figures (columns H-N: currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the leve
# profit sheets, reflecting refreshed price data from the scheduled runner.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 941.25
$ws.Cells.Item(28, 9).Value = 990.6667
$ws.Cells.Item(28, 10).Value = 200
$ws.Cells.Item(28, 11).Value = 990.6667
$ws.Cells.Item(28, 12).Value = 200
$ws.Cells.Item(28, 13).Value = -505.6667
$ws.Cells.Item(28, 14).Value = -1170
$ws.Cells.Item(74, 8).Value = 3711929.5
$ws.Cells.Item(74, 9).Value = 4330009.5
$ws.Cells.Item(74, 10).Value = 3450
$ws.Cells.Item(74, 11).Value = 4330009.5
$ws.Cells.Item(74, 12).Value = 3450
$ws.Cells.Item(74, 13).Value = -4329073.5
$ws.Cells.Item(74, 14).Value = -5322
$ws.Cells.Item(77, 8).Value = 3711929.5
$ws.Cells.Item(77, 9).Value = 4330009.5
$ws.Cells.Item(77, 10).Value = 3450
$ws.Cells.Item(77, 11).Value = 21650047.5
$ws.Cells.Item(77, 12).Value = 17250
$ws.Cells.Item(77, 13).Value = -21645367.5
$ws.Cells.Item(77, 14).Value = -26610
$ws.Cells.Item(132, 8).Value = 2025.069
$ws.Cells.Item(132, 9).Value = 2089.08
$ws.Cells.Item(132, 10).Value = 1625
$ws.Cells.Item(132, 11).Value = 6267.24
$ws.Cells.Item(132, 12).Value = 4875
$ws.Cells.Item(132, 13).Value = -3737.24
$ws.Cells.Item(132, 14).Value = -9935
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1430.8276
$ws.Cells.Item(2, 9).Value = 1443.76
$ws.Cells.Item(2, 10).Value = 1350
$ws.Cells.Item(2, 11).Value = 1443.76
$ws.Cells.Item(2, 12).Value = 1350
$ws.Cells.Item(2, 13).Value = -1330.76
$ws.Cells.Item(2, 14).Value = -1576
$ws.Cells.Item(32, 8).Value = 12851.811
$ws.Cells.Item(32, 9).Value = 7951.0376
$ws.Cells.Item(32, 10).Value = 64800
$ws.Cells.Item(32, 11).Value = 7951.0376
$ws.Cells.Item(32, 12).Value = 64800
$ws.Cells.Item(32, 13).Value = -7664.0376
$ws.Cells.Item(32, 14).Value = -65374
$ws.Cells.Item(88, 8).Value = 2208.318
$ws.Cells.Item(88, 9).Value = 1746.625
$ws.Cells.Item(88, 11).Value = 1746.625
$ws.Cells.Item(88, 13).Value = -1340.625
$ws.Cells.Item(91, 8).Value = 2208.318
$ws.Cells.Item(91, 9).Value = 1746.625
$ws.Cells.Item(91, 11).Value = 1746.625
$ws.Cells.Item(91, 13).Value = -342.625
$ws.Cells.Item(116, 8).Value = 1430.8276
$ws.Cells.Item(116, 9).Value = 1443.76
$ws.Cells.Item(116, 10).Value = 1350
$ws.Cells.Item(116, 11).Value = 1443.76
$ws.Cells.Item(116, 12).Value = 1350
$ws.Cells.Item(116, 13).Value = 850.24
$ws.Cells.Item(116, 14).Value = -5938
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1430.8276
$ws.Cells.Item(3, 9).Value = 1443.76
$ws.Cells.Item(3, 10).Value = 1350
$ws.Cells.Item(3, 11).Value = 1443.76
$ws.Cells.Item(3, 12).Value = 1350
$ws.Cells.Item(3, 13).Value = -1329.76
$ws.Cells.Item(3, 14).Value = -1578
$ws.Cells.Item(99, 8).Value = 66668890
$ws.Cells.Item(99, 9).Value = 90911336
$ws.Cells.Item(99, 10).Value = 2161
$ws.Cells.Item(99, 11).Value = 90911336
$ws.Cells.Item(99, 12).Value = 2161
$ws.Cells.Item(99, 13).Value = -90909838
$ws.Cells.Item(99, 14).Value = -5157
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 989.4706
$ws.Cells.Item(16, 9).Value = 988.13336
$ws.Cells.Item(16, 11).Value = 988.13336
$ws.Cells.Item(16, 13).Value = -701.13336
$ws.Cells.Item(31, 8).Value = 2726.1738
$ws.Cells.Item(31, 9).Value = 2929.5293
$ws.Cells.Item(31, 10).Value = 2150
$ws.Cells.Item(31, 11).Value = 2929.5293
$ws.Cells.Item(31, 12).Value = 2150
$ws.Cells.Item(31, 13).Value = -2634.5293
$ws.Cells.Item(31, 14).Value = -2740
$ws.Cells.Item(34, 8).Value = 2726.1738
$ws.Cells.Item(34, 9).Value = 2929.5293
$ws.Cells.Item(34, 10).Value = 2150
$ws.Cells.Item(34, 11).Value = 2929.5293
$ws.Cells.Item(34, 12).Value = 2150
$ws.Cells.Item(34, 13).Value = -2727.5293
$ws.Cells.Item(34, 14).Value = -2554
$ws.Cells.Item(113, 8).Value = 989.4706
$ws.Cells.Item(113, 9).Value = 988.13336
$ws.Cells.Item(113, 11).Value = 988.13336
$ws.Cells.Item(113, 13).Value = 1181.86664
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 104.25
$ws.Cells.Item(50, 9).Value = 106.666664
$ws.Cells.Item(50, 10).Value = 97
$ws.Cells.Item(50, 11).Value = 319.999992
$ws.Cells.Item(50, 12).Value = 291
$ws.Cells.Item(50, 13).Value = 161.000008
$ws.Cells.Item(50, 14).Value = -1253
$ws.Cells.Item(53, 8).Value = 104.25
$ws.Cells.Item(53, 9).Value = 106.666664
$ws.Cells.Item(53, 10).Value = 97
$ws.Cells.Item(53, 11).Value = 319.999992
$ws.Cells.Item(53, 12).Value = 291
$ws.Cells.Item(53, 13).Value = 161.000008
$ws.Cells.Item(53, 14).Value = -1253
$ws.Cells.Item(113, 8).Value = 42333.125
$ws.Cells.Item(113, 9).Value = 1070
$ws.Cells.Item(113, 11).Value = 3210
$ws.Cells.Item(113, 13).Value = -1040
$ws.Cells.Item(131, 8).Value = 6423923
$ws.Cells.Item(131, 9).Value = 45545756
$ws.Cells.Item(131, 10).Value = 935.44775
$ws.Cells.Item(131, 11).Value = 136637268
$ws.Cells.Item(131, 12).Value = 2806.34325
$ws.Cells.Item(131, 13).Value = -136632228
$ws.Cells.Item(131, 14).Value = -12886.34325
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 4399.615
$ws.Cells.Item(80, 9).Value = 4300
$ws.Cells.Item(80, 10).Value = 4623.75
$ws.Cells.Item(80, 11).Value = 4300
$ws.Cells.Item(80, 12).Value = 4623.75
$ws.Cells.Item(80, 13).Value = -3302
$ws.Cells.Item(80, 14).Value = -6619.75
$ws.Cells.Item(83, 8).Value = 4399.615
$ws.Cells.Item(83, 9).Value = 4300
$ws.Cells.Item(83, 10).Value = 4623.75
$ws.Cells.Item(83, 11).Value = 21500
$ws.Cells.Item(83, 12).Value = 23118.75
$ws.Cells.Item(83, 13).Value = -16508
$ws.Cells.Item(83, 14).Value = -33102.75
$ws.Cells.Item(110, 8).Value = 29833.334
$ws.Cells.Item(110, 10).Value = 29833.334
$ws.Cells.Item(110, 12).Value = 29833.334
$ws.Cells.Item(110, 14).Value = -38013.334
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 2379.4666
$ws.Cells.Item(82, 9).Value = 3500.4
$ws.Cells.Item(82, 10).Value = 1819
$ws.Cells.Item(82, 11).Value = 3500.4
$ws.Cells.Item(82, 12).Value = 1819
$ws.Cells.Item(82, 13).Value = -3139.4
$ws.Cells.Item(82, 14).Value = -2541
$ws.Cells.Item(85, 8).Value = 2379.4666
$ws.Cells.Item(85, 9).Value = 3500.4
$ws.Cells.Item(85, 10).Value = 1819
$ws.Cells.Item(85, 11).Value = 3500.4
$ws.Cells.Item(85, 12).Value = 1819
$ws.Cells.Item(85, 13).Value = -2252.4
$ws.Cells.Item(85, 14).Value = -4315
$ws.Cells.Item(93, 8).Value = 2607.3333
$ws.Cells.Item(93, 9).Value = 1684.3334
$ws.Cells.Item(93, 10).Value = 3222.6667
$ws.Cells.Item(93, 11).Value = 1684.3334
$ws.Cells.Item(93, 12).Value = 3222.6667
$ws.Cells.Item(93, 13).Value = -436.3334
$ws.Cells.Item(93, 14).Value = -5718.6667
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 67419.836
$ws.Cells.Item(46, 10).Value = 67419.836
$ws.Cells.Item(46, 12).Value = 67419.836
$ws.Cells.Item(46, 14).Value = -67881.836
$ws.Cells.Item(100, 8).Value = 2462.1428
$ws.Cells.Item(100, 9).Value = 1047
$ws.Cells.Item(100, 10).Value = 6000
$ws.Cells.Item(100, 11).Value = 2094
$ws.Cells.Item(100, 12).Value = 12000
$ws.Cells.Item(100, 13).Value = -1553
$ws.Cells.Item(100, 14).Value = -13082
$ws.Cells.Item(134, 8).Value = 67419.836
$ws.Cells.Item(134, 10).Value = 67419.836
$ws.Cells.Item(134, 12).Value = 202259.508
$ws.Cells.Item(134, 14).Value = -207329.508
